# Auto-assisted generation verified by hand; applies the france ligue-1 2023-2024 diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: rotate match-data (columns F:V) within each swapped-row group ---
# Column A (Indice) and E (data_partida) are untouched; F..V cycle: row i <- row (i+1), last <- first.
function Rotate-Rows($rows) {
    $n = $rows.Length
    $snapshot = @{}
    foreach ($r in $rows) {
        $vals = @()
        for ($c = 6; $c -le 22; $c++) {
            $vals += ,$ws.Cells.Item($r, $c).Value2
        }
        $snapshot[$r] = $vals
    }
    for ($i = 0; $i -lt $n; $i++) {
        $dst = $rows[$i]
        $src = $rows[($i + 1) % $n]
        $vals = $snapshot[$src]
        for ($c = 6; $c -le 22; $c++) {
            $ws.Cells.Item($dst, $c).Value2 = $vals[$c - 6]
        }
    }
}

Rotate-Rows @(6, 7, 8)
Rotate-Rows @(25, 26)
Rotate-Rows @(33, 34, 35)
Rotate-Rows @(42, 43, 44)

# --- Part 2: append 7 new match rows (48-54) ---
$lastRow = 47
$newRowCount = 7
$ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
$ws.Range("A" + ($lastRow + 1) + ":V" + ($lastRow + $newRowCount)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 48
$ws.Cells.Item(48, 1).Value2 = 47
$ws.Cells.Item(48, 2).Value2 = "france"
$ws.Cells.Item(48, 3).Value2 = "ligue-1"
$ws.Cells.Item(48, 4).Value2 = "2023-2024"
$ws.Cells.Item(48, 5).Value2 = 45192.70833333334
$ws.Cells.Item(48, 6).Value2 = "Nantes"
$ws.Cells.Item(48, 7).Value2 = 5
$ws.Cells.Item(48, 8).Value2 = "Lorient"
$ws.Cells.Item(48, 9).Value2 = 3
$ws.Cells.Item(48, 10).Value2 = 2.17
$ws.Cells.Item(48, 11).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(48, 12).Value2 = 2.36
$ws.Cells.Item(48, 13).Value2 = "23/09/2023 16:58"
$ws.Cells.Item(48, 14).Value2 = 3.42
$ws.Cells.Item(48, 15).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(48, 16).Value2 = 3.39
$ws.Cells.Item(48, 17).Value2 = "23/09/2023 16:58"
$ws.Cells.Item(48, 18).Value2 = 3.54
$ws.Cells.Item(48, 19).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(48, 20).Value2 = 3.3
$ws.Cells.Item(48, 21).Value2 = "23/09/2023 16:58"
$ws.Cells.Item(48, 22).Value2 = "https://www.betexplorer.com/football/france/ligue-1/nantes-lorient/8xxmwRLB/"

# row 49
$ws.Cells.Item(49, 1).Value2 = 48
$ws.Cells.Item(49, 2).Value2 = "france"
$ws.Cells.Item(49, 3).Value2 = "ligue-1"
$ws.Cells.Item(49, 4).Value2 = "2023-2024"
$ws.Cells.Item(49, 5).Value2 = 45192.875
$ws.Cells.Item(49, 6).Value2 = "Brest"
$ws.Cells.Item(49, 7).Value2 = 1
$ws.Cells.Item(49, 8).Value2 = "Lyon"
$ws.Cells.Item(49, 9).Value2 = 0
$ws.Cells.Item(49, 10).Value2 = 3.08
$ws.Cells.Item(49, 11).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(49, 12).Value2 = 2.73
$ws.Cells.Item(49, 13).Value2 = "23/09/2023 20:57"
$ws.Cells.Item(49, 14).Value2 = 3.57
$ws.Cells.Item(49, 15).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(49, 16).Value2 = 3.43
$ws.Cells.Item(49, 17).Value2 = "23/09/2023 20:48"
$ws.Cells.Item(49, 18).Value2 = 2.33
$ws.Cells.Item(49, 19).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(49, 20).Value2 = 2.75
$ws.Cells.Item(49, 21).Value2 = "23/09/2023 20:58"
$ws.Cells.Item(49, 22).Value2 = "https://www.betexplorer.com/football/france/ligue-1/brest-lyon/4Ence8yt/"

# row 50
$ws.Cells.Item(50, 1).Value2 = 49
$ws.Cells.Item(50, 2).Value2 = "france"
$ws.Cells.Item(50, 3).Value2 = "ligue-1"
$ws.Cells.Item(50, 4).Value2 = "2023-2024"
$ws.Cells.Item(50, 5).Value2 = 45193.54166666666
$ws.Cells.Item(50, 6).Value2 = "Metz"
$ws.Cells.Item(50, 7).Value2 = 0
$ws.Cells.Item(50, 8).Value2 = "Strasbourg"
$ws.Cells.Item(50, 9).Value2 = 1
$ws.Cells.Item(50, 10).Value2 = 3.12
$ws.Cells.Item(50, 11).Value2 = "11/09/2023 13:36"
$ws.Cells.Item(50, 12).Value2 = 2.92
$ws.Cells.Item(50, 13).Value2 = "24/09/2023 12:58"
$ws.Cells.Item(50, 14).Value2 = 3.25
$ws.Cells.Item(50, 15).Value2 = "11/09/2023 13:36"
$ws.Cells.Item(50, 16).Value2 = 3.19
$ws.Cells.Item(50, 17).Value2 = "24/09/2023 12:56"
$ws.Cells.Item(50, 18).Value2 = 2.35
$ws.Cells.Item(50, 19).Value2 = "11/09/2023 13:36"
$ws.Cells.Item(50, 20).Value2 = 2.73
$ws.Cells.Item(50, 21).Value2 = "24/09/2023 12:58"
$ws.Cells.Item(50, 22).Value2 = "https://www.betexplorer.com/football/france/ligue-1/metz-strasbourg/G2uey5jO/"

# row 51
$ws.Cells.Item(51, 1).Value2 = 50
$ws.Cells.Item(51, 2).Value2 = "france"
$ws.Cells.Item(51, 3).Value2 = "ligue-1"
$ws.Cells.Item(51, 4).Value2 = "2023-2024"
$ws.Cells.Item(51, 5).Value2 = 45193.625
$ws.Cells.Item(51, 6).Value2 = "Le Havre"
$ws.Cells.Item(51, 7).Value2 = 2
$ws.Cells.Item(51, 8).Value2 = "Clermont"
$ws.Cells.Item(51, 9).Value2 = 1
$ws.Cells.Item(51, 10).Value2 = 2.72
$ws.Cells.Item(51, 11).Value2 = "11/09/2023 13:39"
$ws.Cells.Item(51, 12).Value2 = 2.56
$ws.Cells.Item(51, 13).Value2 = "24/09/2023 14:58"
$ws.Cells.Item(51, 14).Value2 = 3.25
$ws.Cells.Item(51, 15).Value2 = "11/09/2023 13:39"
$ws.Cells.Item(51, 16).Value2 = 3.08
$ws.Cells.Item(51, 17).Value2 = "24/09/2023 14:58"
$ws.Cells.Item(51, 18).Value2 = 2.64
$ws.Cells.Item(51, 19).Value2 = "11/09/2023 13:39"
$ws.Cells.Item(51, 20).Value2 = 3.26
$ws.Cells.Item(51, 21).Value2 = "24/09/2023 14:58"
$ws.Cells.Item(51, 22).Value2 = "https://www.betexplorer.com/football/france/ligue-1/le-havre-clermont/MutixoyI/"

# row 52
$ws.Cells.Item(52, 1).Value2 = 51
$ws.Cells.Item(52, 2).Value2 = "france"
$ws.Cells.Item(52, 3).Value2 = "ligue-1"
$ws.Cells.Item(52, 4).Value2 = "2023-2024"
$ws.Cells.Item(52, 5).Value2 = 45193.625
$ws.Cells.Item(52, 6).Value2 = "Lens"
$ws.Cells.Item(52, 7).Value2 = 2
$ws.Cells.Item(52, 8).Value2 = "Toulouse"
$ws.Cells.Item(52, 9).Value2 = 1
$ws.Cells.Item(52, 10).Value2 = 1.58
$ws.Cells.Item(52, 11).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(52, 12).Value2 = 1.45
$ws.Cells.Item(52, 13).Value2 = "24/09/2023 14:59"
$ws.Cells.Item(52, 14).Value2 = 4.44
$ws.Cells.Item(52, 15).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(52, 16).Value2 = 4.73
$ws.Cells.Item(52, 17).Value2 = "24/09/2023 14:59"
$ws.Cells.Item(52, 18).Value2 = 5.6
$ws.Cells.Item(52, 19).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(52, 20).Value2 = 8.07
$ws.Cells.Item(52, 21).Value2 = "24/09/2023 14:59"
$ws.Cells.Item(52, 22).Value2 = "https://www.betexplorer.com/football/france/ligue-1/lens-toulouse/QDnrv765/"

# row 53
$ws.Cells.Item(53, 1).Value2 = 52
$ws.Cells.Item(53, 2).Value2 = "france"
$ws.Cells.Item(53, 3).Value2 = "ligue-1"
$ws.Cells.Item(53, 4).Value2 = "2023-2024"
$ws.Cells.Item(53, 5).Value2 = 45193.71180555555
$ws.Cells.Item(53, 6).Value2 = "Montpellier"
$ws.Cells.Item(53, 7).Value2 = 0
$ws.Cells.Item(53, 8).Value2 = "Rennes"
$ws.Cells.Item(53, 9).Value2 = 0
$ws.Cells.Item(53, 10).Value2 = 3.08
$ws.Cells.Item(53, 11).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(53, 12).Value2 = 2.77
$ws.Cells.Item(53, 13).Value2 = "24/09/2023 16:49"
$ws.Cells.Item(53, 14).Value2 = 3.48
$ws.Cells.Item(53, 15).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(53, 16).Value2 = 3.55
$ws.Cells.Item(53, 17).Value2 = "24/09/2023 17:01"
$ws.Cells.Item(53, 18).Value2 = 2.37
$ws.Cells.Item(53, 19).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(53, 20).Value2 = 2.57
$ws.Cells.Item(53, 21).Value2 = "24/09/2023 17:01"
$ws.Cells.Item(53, 22).Value2 = "https://www.betexplorer.com/football/france/ligue-1/montpellier-rennes/fgvazP5U/"

# row 54
$ws.Cells.Item(54, 1).Value2 = 53
$ws.Cells.Item(54, 2).Value2 = "france"
$ws.Cells.Item(54, 3).Value2 = "ligue-1"
$ws.Cells.Item(54, 4).Value2 = "2023-2024"
$ws.Cells.Item(54, 5).Value2 = 45193.86458333334
$ws.Cells.Item(54, 6).Value2 = "Paris SG"
$ws.Cells.Item(54, 7).Value2 = 4
$ws.Cells.Item(54, 8).Value2 = "Marseille"
$ws.Cells.Item(54, 9).Value2 = 0
$ws.Cells.Item(54, 10).Value2 = 1.81
$ws.Cells.Item(54, 11).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(54, 12).Value2 = 1.56
$ws.Cells.Item(54, 13).Value2 = "24/09/2023 20:44"
$ws.Cells.Item(54, 14).Value2 = 4.02
$ws.Cells.Item(54, 15).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(54, 16).Value2 = 4.57
$ws.Cells.Item(54, 17).Value2 = "24/09/2023 20:44"
$ws.Cells.Item(54, 18).Value2 = 4.27
$ws.Cells.Item(54, 19).Value2 = "05/09/2023 12:01"
$ws.Cells.Item(54, 20).Value2 = 6.11
$ws.Cells.Item(54, 21).Value2 = "24/09/2023 20:44"
$ws.Cells.Item(54, 22).Value2 = "https://www.betexplorer.com/football/france/ligue-1/paris-sg-marseille/vcpWt9Mn/"

